$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price/volume snapshot (Price = column D,
# Volume(1h) = column E). Rows without a Price change only touch column E.
# For Price values that read as plain numbers (e.g. "1.00", "7.82"), the
# cell's NumberFormat is briefly switched to Text ("@") before the write so
# the value is stored verbatim as a string instead of being parsed into a
# number (which would silently drop the trailing/formatting zeros); the
# style is then reset to "Normal" so no visible formatting change remains.

$ws.Range("D2").Value = '51.691.49'
$ws.Range("E2").Value = '  +0.25%  '

$ws.Range("D3").Value = '3.082.18'
$ws.Range("E3").Value = '  +3.16%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '388.63'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.64%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.59'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.01%  '

$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.588'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.01%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.07'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.91%  '

$ws.Range("E11").Value = '  +0.23%  '

$ws.Range("E12").Value = '  +0.63%  '

$ws.Range("D13").Value = '3.567.00'
$ws.Range("E13").Value = '  +3.59%  '

$ws.Range("E14").Value = '  +1.74%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.82'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.15%  '

$ws.Range("D16").Value = '3.078.11'
$ws.Range("E16").Value = '  +3.11%  '

$ws.Range("E17").Value = '  -2.02%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.75'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.89%  '

$ws.Range("D19").Value = '51.787.20'
$ws.Range("E19").Value = '  +0.49%  '

$ws.Range("E20").Value = '  +3.01%  '

$ws.Range("E21").Value = '  -1.31%  '

$ws.Range("D22").Value = '0.0₃0970'
$ws.Range("E22").Value = '  +0.60%  '

$ws.Range("E23").Value = '  -0.16%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '268.96'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.53%  '

$ws.Range("E25").Value = '  -1.77%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.23'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.16%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '27.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.34%  '

$ws.Range("E28").Value = '  -0.31%  '

$ws.Range("E29").Value = '  +1.22%  '

$ws.Range("E30").Value = '  +0.05%  '

$ws.Range("E31").Value = '  -1.29%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '10.29'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.41%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '34.98'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.55%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.08'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.30%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '50.15'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.80%  '

$ws.Range("E36").Value = '  +1.89%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.03%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.34'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.79%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.296'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +9.32%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.89'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.25%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '16.97'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.04%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.57'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.32%  '

$ws.Range("E43").Value = '  -0.84%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '125.98'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.51%  '

$ws.Range("E45").Value = '  -0.78%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '21.98'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.14%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.10'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.98%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.45'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.37%  '

$ws.Range("D49").Value = '2.040.65'
$ws.Range("E49").Value = '  +0.55%  '

$ws.Range("D50").Value = '3.381.90'
$ws.Range("E50").Value = '  +3.20%  '

$ws.Range("E51").Value = '  +6.84%  '
